
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Cebollín" at
# "Feria Lagunitas de Puerto Montt". It belongs chronologically at the top
# of this market/category block (row 466), so insert a fresh row there and
# push the existing 466:553 block (and all data below it) down by one row.
$ws.Rows("466:466").Insert()

$ws.Range("A466").Value = 4
$ws.Range("B466").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C466").Value = "Los Lagos"
$ws.Range("D466").Value = 45258
$ws.Range("E466").Value = 10
$ws.Range("F466").Value = 100112037
$ws.Range("G466").Value = "Cebollín"
$ws.Range("H466").Value = "Sin especificar"
$ws.Range("I466").Value = "Primera"
$ws.Range("J466").Value = 180
$ws.Range("K466").Value = 6000
$ws.Range("L466").Value = 6500
$ws.Range("M466").Value = 6250
$ws.Range("N466").Value = "`$/paquete 36 unidades"
$ws.Range("O466").Value = "Región Metropolitana"
$ws.Range("P466").Value = 174
$ws.Range("Q466").Value = 36
$ws.Range("R466").Value = "Hortaliza"
